# daily auto push: 2025-10-08 06:41 UTC
# Append the new day's row (row 78) to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 78

# Column A holds a date-looking string ("2025/10/08") that must stay plain
# text, exactly like the existing rows (e.g. A76/A77), not get auto-converted
# into a date serial number. Force text formatting before assigning it, then
# clear the formatting override afterwards so the cell keeps the default
# (unstyled) look of its neighbours.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/08"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "水"
$ws.Cells.Item($row, 3).Value = 14
$ws.Cells.Item($row, 4).Value = 115
